# Auto update Excel log
#
# Appends new sensor-log rows scraped from the monitoring system to the
# PIR, Humidity, Proximity and Camera sheets.
#
# NOTE: Assigning a literal ".Value" that *looks* like a date ("2026-01-30"),
# a time ("16:09:16") or a percentage ("86.1%") makes Excel silently
# reinterpret it as a real date/time/number and stamp the cell with a
# non-default number format. To keep the new cells as plain text (matching
# every other cell already in these logs) each cell is written using the
# standard Excel workaround:
#   1. Force the cell's number format to Text ("@") so the literal string
#      is stored verbatim instead of being parsed.
#   2. Assign the literal value.
#   3. Reset the cell style back to "Normal" so the cell keeps the same
#      (default/general) appearance as the rest of the sheet.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# PIR sheet: append rows 71-83 (Bathroom / No Motion / Inactive)
# ----------------------------------------------------------------------------
$wsPIR = $wb.Worksheets.Item("PIR")
$pirLines = @(
"2026-01-30|16:09:16|16:00|Bathroom|No Motion|Inactive",
"2026-01-30|16:09:17|16:00|Bathroom|No Motion|Inactive",
"2026-01-30|16:09:22|16:00|Bathroom|No Motion|Inactive",
"2026-01-30|16:09:27|16:00|Bathroom|No Motion|Inactive",
"2026-01-30|16:09:32|16:00|Bathroom|No Motion|Inactive",
"2026-01-30|16:09:37|16:00|Bathroom|No Motion|Inactive",
"2026-01-30|16:09:42|16:00|Bathroom|No Motion|Inactive",
"2026-01-30|16:09:47|16:00|Bathroom|No Motion|Inactive",
"2026-01-30|16:09:52|16:00|Bathroom|No Motion|Inactive",
"2026-01-30|16:09:57|16:00|Bathroom|No Motion|Inactive",
"2026-01-30|16:10:02|16:00|Bathroom|No Motion|Inactive",
"2026-01-30|16:10:07|16:00|Bathroom|No Motion|Inactive",
"2026-01-30|16:10:12|16:00|Bathroom|No Motion|Inactive"
)
$r = 71
foreach ($line in $pirLines) {
    $parts = $line.Split("|")
    for ($col = 1; $col -le 6; $col++) {
        $cell = $wsPIR.Cells.Item($r, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $parts[$col - 1]
        $cell.Style = "Normal"
    }
    $r = $r + 1
}

# ----------------------------------------------------------------------------
# Humidity sheet: append rows 56-67 (Bathroom / xx.x% / Active)
# ----------------------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humidityLines = @(
"2026-01-30|16:09:16|16:00|Bathroom|86.1%|Active",
"2026-01-30|16:09:17|16:00|Bathroom|87.6%|Active",
"2026-01-30|16:09:27|16:00|Bathroom|87.6%|Active",
"2026-01-30|16:09:32|16:00|Bathroom|86.6%|Active",
"2026-01-30|16:09:37|16:00|Bathroom|87.5%|Active",
"2026-01-30|16:09:42|16:00|Bathroom|86.6%|Active",
"2026-01-30|16:09:47|16:00|Bathroom|87.5%|Active",
"2026-01-30|16:09:52|16:00|Bathroom|87.6%|Active",
"2026-01-30|16:09:57|16:00|Bathroom|87.5%|Active",
"2026-01-30|16:10:02|16:00|Bathroom|86.6%|Active",
"2026-01-30|16:10:07|16:00|Bathroom|87.6%|Active",
"2026-01-30|16:10:12|16:00|Bathroom|86.6%|Active"
)
$r = 56
foreach ($line in $humidityLines) {
    $parts = $line.Split("|")
    for ($col = 1; $col -le 6; $col++) {
        $cell = $wsHumidity.Cells.Item($r, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $parts[$col - 1]
        $cell.Style = "Normal"
    }
    $r = $r + 1
}

# ----------------------------------------------------------------------------
# Proximity sheet: append row 16 (Living Room Main Door / ENTER)
# ----------------------------------------------------------------------------
$wsProximity = $wb.Worksheets.Item("Proximity")
$parts = "2026-01-30|16:09:34|16:00|Living Room Main Door|ENTER|User ENTERED Living Room Main Door".Split("|")
for ($col = 1; $col -le 6; $col++) {
    $cell = $wsProximity.Cells.Item(16, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $parts[$col - 1]
    $cell.Style = "Normal"
}

# ----------------------------------------------------------------------------
# Camera sheet: append row 16 (Living Room Main Door / Image Captured (ENTER))
# ----------------------------------------------------------------------------
$wsCamera = $wb.Worksheets.Item("Camera")
$parts = "2026-01-30|16:09:33|16:00|Living Room Main Door|Image Captured (ENTER)|Active".Split("|")
for ($col = 1; $col -le 6; $col++) {
    $cell = $wsCamera.Cells.Item(16, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $parts[$col - 1]
    $cell.Style = "Normal"
}
